# Update "想去人数" (F column) counters across the workbook sheets.
# Sheet order in the workbook: 1=展览, 2=演出, 3=本地生活, 4=全部类型

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 2716
$ws1.Range("F5").Value = 216
$ws1.Range("F7").Value = 1295
$ws1.Range("F8").Value = 611
$ws1.Range("F13").Value = 5984
$ws1.Range("F16").Value = 1858
$ws1.Range("F17").Value = 4449
$ws1.Range("F18").Value = 453
$ws1.Range("F21").Value = 5187
$ws1.Range("F22").Value = 6730
$ws1.Range("F26").Value = 3902
$ws1.Range("F27").Value = 526
$ws1.Range("F29").Value = 210
$ws1.Range("F32").Value = 1461
$ws1.Range("F34").Value = 630
$ws1.Range("F35").Value = 1650
$ws1.Range("F36").Value = 221
$ws1.Range("F37").Value = 1819
$ws1.Range("F39").Value = 1196
$ws1.Range("F41").Value = 656
$ws1.Range("F43").Value = 3568
$ws1.Range("F45").Value = 326
$ws1.Range("F49").Value = 3922

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 1241

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 4203

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 4203
$ws4.Range("F3").Value = 2716
$ws4.Range("F6").Value = 1241
$ws4.Range("F8").Value = 216
$ws4.Range("F11").Value = 1295
$ws4.Range("F13").Value = 611
$ws4.Range("F18").Value = 1858
$ws4.Range("F19").Value = 4449
$ws4.Range("F20").Value = 5187
$ws4.Range("F21").Value = 5187
$ws4.Range("F25").Value = 3902
$ws4.Range("F26").Value = 526
$ws4.Range("F27").Value = 210
$ws4.Range("F30").Value = 1461
$ws4.Range("F32").Value = 630
$ws4.Range("F33").Value = 1650
$ws4.Range("F34").Value = 221
$ws4.Range("F35").Value = 1819
$ws4.Range("F39").Value = 656
$ws4.Range("F43").Value = 3568
$ws4.Range("F46").Value = 326
$ws4.Range("F50").Value = 3922
